# Adição de novos arquivos
# - Preenche o número de responsável que faltava (B11)
# - Alarga a coluna A para caber nomes maiores (sem autofit/bestFit)
# - Ajusta a visualização da planilha (rolagem/seleção)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# B11 estava vazio; adiciona o número de telefone do responsável.
$ws.Range("B11").Value = 5543984250977

# Coluna A: largura explícita de 46 (remove o ajuste automático "bestFit").
# A largura armazenada no XML fica 5/6 acima do valor de ColumnWidth nesta
# engine, então compensamos para gravar exatamente width="46".
$ws.Columns("A").ColumnWidth = 46 - 5/6

# Rola a planilha e seleciona a célula F14, como deixado pelo autor.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F14").Select() | Out-Null
